$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data of rows 4-5 with rows 6-7 (Belegnummer, Belegdatum, Lieferant, Profitcenter, Kostenstelle)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "05.Jan"
$ws.Range("C4").Value = "A"
$ws.Range("D4").Value = "P1"
$ws.Range("E4").Value = "K2"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "07.Jan"
$ws.Range("C5").Value = "A"
$ws.Range("D5").Value = "P1"
$ws.Range("E5").Value = "K2"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "09.Jan"
$ws.Range("C6").Value = "B"
$ws.Range("D6").Value = "P2"
$ws.Range("E6").Value = "K1"

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "11.Jan"
$ws.Range("C7").Value = "A"
$ws.Range("D7").Value = "P2"
$ws.Range("E7").Value = "K1"
